$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vtn"
$ws.Range("C2").Value = "Itgav"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 6.597131
$ws.Range("H2").Value = 19.791393
$ws.Range("I2").Value = 0.1209543635982448
$ws.Range("J2").Value = 0.1209543635982448
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 15.35884066666667
$ws.Range("N2").Value = 46.076522
$ws.Range("O2").Value = 0.1012042817263867
$ws.Range("P2").Value = 0.1012042817263867
$ws.Range("Q2").Value = 101.3242838861273
$ws.Range("R2").Value = 911.9185549751459
$ws.Range("S2").Value = 0.01224109948963257
$ws.Range("T2").Value = 0.01224109948963257

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vtn"
$ws.Range("C3").Value = "Itgav"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 6.597131
$ws.Range("H3").Value = 19.791393
$ws.Range("I3").Value = 0.1209543635982448
$ws.Range("J3").Value = 0.1209543635982448
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 50.59256466666667
$ws.Range("N3").Value = 151.777694
$ws.Range("O3").Value = 0.3333704853712116
$ws.Range("P3").Value = 0.3333704853712116
$ws.Range("Q3").Value = 333.7657767319714
$ws.Range("R3").Value = 3003.891990587742
$ws.Range("S3").Value = 0.04032261490051286
$ws.Range("T3").Value = 0.04032261490051286

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vtn"
$ws.Range("C4").Value = "Itgav"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 6.597131
$ws.Range("H4").Value = 19.791393
$ws.Range("I4").Value = 0.1209543635982448
$ws.Range("J4").Value = 0.1209543635982448
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 60.37715666666667
$ws.Range("N4").Value = 181.13147
$ws.Range("O4").Value = 0.397844271305776
$ws.Range("P4").Value = 0.397844271305776
$ws.Range("Q4").Value = 398.3160119375233
$ws.Range("R4").Value = 3584.84410743771
$ws.Range("S4").Value = 0.04812100064699756
$ws.Range("T4").Value = 0.04812100064699756

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Vtn"
$ws.Range("C5").Value = "Itgav"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 6.597131
$ws.Range("H5").Value = 19.791393
$ws.Range("I5").Value = 0.1209543635982448
$ws.Range("J5").Value = 0.1209543635982448
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 25.43221733333333
$ws.Range("N5").Value = 76.296652
$ws.Range("O5").Value = 0.1675809615966257
$ws.Range("P5").Value = 0.1675809615966258
$ws.Range("Q5").Value = 167.7796693684706
$ws.Range("R5").Value = 1510.017024316236
$ws.Range("S5").Value = 0.02026964856110176
$ws.Range("T5").Value = 0.02026964856110177

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vtn"
$ws.Range("C6").Value = "Itgav"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 21.05317333333333
$ws.Range("H6").Value = 63.15952
$ws.Range("I6").Value = 0.3859970617919927
$ws.Range("J6").Value = 0.3859970617919927
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 15.35884066666667
$ws.Range("N6").Value = 46.076522
$ws.Range("O6").Value = 0.1012042817263867
$ws.Range("P6").Value = 0.1012042817263867
$ws.Range("Q6").Value = 323.3523347543822
$ws.Range("R6").Value = 2910.17101278944
$ws.Range("S6").Value = 0.03906455538715431
$ws.Range("T6").Value = 0.03906455538715432

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vtn"
$ws.Range("C7").Value = "Itgav"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 21.05317333333333
$ws.Range("H7").Value = 63.15952
$ws.Range("I7").Value = 0.3859970617919927
$ws.Range("J7").Value = 0.3859970617919927
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 50.59256466666667
$ws.Range("N7").Value = 151.777694
$ws.Range("O7").Value = 0.3333704853712116
$ws.Range("P7").Value = 0.3333704853712116
$ws.Range("Q7").Value = 1065.134033305209
$ws.Range("R7").Value = 9586.20629974688
$ws.Range("S7").Value = 0.1286800278414582
$ws.Range("T7").Value = 0.1286800278414582

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Vtn"
$ws.Range("C8").Value = "Itgav"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 21.05317333333333
$ws.Range("H8").Value = 63.15952
$ws.Range("I8").Value = 0.3859970617919927
$ws.Range("J8").Value = 0.3859970617919927
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 60.37715666666667
$ws.Range("N8").Value = 181.13147
$ws.Range("O8").Value = 0.397844271305776
$ws.Range("P8").Value = 0.397844271305776
$ws.Range("Q8").Value = 1271.130744677156
$ws.Range("R8").Value = 11440.1767020944
$ws.Range("S8").Value = 0.1535667197748059
$ws.Range("T8").Value = 0.1535667197748059

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Vtn"
$ws.Range("C9").Value = "Itgav"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 21.05317333333333
$ws.Range("H9").Value = 63.15952
$ws.Range("I9").Value = 0.3859970617919927
$ws.Range("J9").Value = 0.3859970617919927
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 25.43221733333333
$ws.Range("N9").Value = 76.296652
$ws.Range("O9").Value = 0.1675809615966257
$ws.Range("P9").Value = 0.1675809615966258
$ws.Range("Q9").Value = 535.4288797696711
$ws.Range("R9").Value = 4818.85991792704
$ws.Range("S9").Value = 0.06468575878857431
$ws.Range("T9").Value = 0.06468575878857431

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Vtn"
$ws.Range("C10").Value = "Itgav"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 26.89201066666667
$ws.Range("H10").Value = 80.676032
$ws.Range("I10").Value = 0.4930485746097625
$ws.Range("J10").Value = 0.4930485746097625
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 15.35884066666667
$ws.Range("N10").Value = 46.076522
$ws.Range("O10").Value = 0.1012042817263867
$ws.Range("P10").Value = 0.1012042817263867
$ws.Range("Q10").Value = 413.0301070356338
$ws.Range("R10").Value = 3717.270963320704
$ws.Range("S10").Value = 0.04989862684959977
$ws.Range("T10").Value = 0.04989862684959978

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Vtn"
$ws.Range("C11").Value = "Itgav"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 26.89201066666667
$ws.Range("H11").Value = 80.676032
$ws.Range("I11").Value = 0.4930485746097625
$ws.Range("J11").Value = 0.4930485746097625
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 50.59256466666667
$ws.Range("N11").Value = 151.777694
$ws.Range("O11").Value = 0.3333704853712116
$ws.Range("P11").Value = 0.3333704853712116
$ws.Range("Q11").Value = 1360.535788670023
$ws.Range("R11").Value = 12244.82209803021
$ws.Range("S11").Value = 0.1643678426292405
$ws.Range("T11").Value = 0.1643678426292405

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Vtn"
$ws.Range("C12").Value = "Itgav"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 26.89201066666667
$ws.Range("H12").Value = 80.676032
$ws.Range("I12").Value = 0.4930485746097625
$ws.Range("J12").Value = 0.4930485746097625
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 60.37715666666667
$ws.Range("N12").Value = 181.13147
$ws.Range("O12").Value = 0.397844271305776
$ws.Range("P12").Value = 0.397844271305776
$ws.Range("Q12").Value = 1623.663141103005
$ws.Range("R12").Value = 14612.96826992704
$ws.Range("S12").Value = 0.1961565508839725
$ws.Range("T12").Value = 0.1961565508839725

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Vtn"
$ws.Range("C13").Value = "Itgav"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 26.89201066666667
$ws.Range("H13").Value = 80.676032
$ws.Range("I13").Value = 0.4930485746097625
$ws.Range("J13").Value = 0.4930485746097625
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 25.43221733333333
$ws.Range("N13").Value = 76.296652
$ws.Range("O13").Value = 0.1675809615966257
$ws.Range("P13").Value = 0.1675809615966258
$ws.Range("Q13").Value = 683.9234598049849
$ws.Range("R13").Value = 6155.311138244864
$ws.Range("S13").Value = 0.08262555424694967
$ws.Range("T13").Value = 0.08262555424694969
